# Reorder the player roster rows on the active sheet.
# The header (row 1) and the overall set of (Player, Position, Team) rows
# are unchanged -- only the row order for rows 2-19 is rearranged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final desired row order (row number -> Player, Position, Team)
$rows = @{
    2  = @("Austin Reaves", "PG,SG", "Los Angeles Lakers")
    3  = @("Jrue Holiday", "PG,SG", "Boston Celtics")
    4  = @("Stephen Curry", "PG,SG", "Golden State Warriors")
    5  = @("Darius Garland", "PG", "Cleveland Cavaliers")
    6  = @("OG Anunoby", "SF,PF", "New York Knicks")
    7  = @("Daniel Gafford", "PF,C", "Dallas Mavericks")
    8  = @("Karl-Anthony Towns", "PF,C", "New York Knicks")
    9  = @("Kevin Durant", "SF,PF", "Phoenix Suns")
    10 = @("Jarrett Allen", "C", "Cleveland Cavaliers")
    11 = @("Jalen Duren", "C", "Detroit Pistons")
    12 = @("Trey Murphy III", "SF,PF", "New Orleans Pelicans")
    13 = @("Mark Williams", "C", "Charlotte Hornets")
    14 = @("Tyrese Haliburton", "PG,SG", "Indiana Pacers")
    15 = @("Keegan Murray", "SF,PF", "Sacramento Kings")
    16 = @("Tyrese Maxey", "PG,SG", "Philadelphia 76ers")
    17 = @("Franz Wagner", "SF,PF", "Orlando Magic")
    18 = @("Dereck Lively II", "C", "Dallas Mavericks")
    19 = @("Jalen Johnson", "SF,PF", "Atlanta Hawks")
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
}
